$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated metrics for the last row (ano/ano_obj = 2025)
$ws.Range("C8").Value = 1160
$ws.Range("D8").Value = 192
$ws.Range("E8").Value = 968
$ws.Range("F8").Value = 7.875307629204266
$ws.Range("G8").Value = 83.44827586206897
$ws.Range("H8").Value = 16.55172413793104
